$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price / Volume(1h) updates for most rows ---
$ws.Range("D2").Value = '47.928.87'
$ws.Range("E2").Value = '  +1.42%  '
$ws.Range("D3").Value = '2.502.07'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '323.57'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").Value = '109.63'
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("D7").Value = '0.525'
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.554'
$ws.Range("E9").Value = '  +2.13%  '
$ws.Range("E10").Value = '  +6.25%  '
$ws.Range("D11").Value = '0.0816'
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '0.124'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").Value = '18.78'
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = '7.24'
$ws.Range("E14").Value = '  +0.51%  '
$ws.Range("D15").Value = '2.892.66'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("D16").Value = '2.504.53'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = '47.831.35'
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("D19").Value = '13.22'
$ws.Range("E19").Value = '  +1.81%  '
$ws.Range("E20").Value = '  -0.53%  '

# --- Row 21/22 swap: ShibaInu <-> ImmutableX ---
$sub3 = [string][char]0x2083
$ws.Range("B21").Value = 'ImmutableX'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D21").Value = '2.79'
$ws.Range("E21").Value = '  +13.50%  '

$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").Value = "0.0" + $sub3 + "0944"
$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("D23").Value = '70.85'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").Value = '248.18'
$ws.Range("E24").Value = '  -0.91%  '
$ws.Range("D25").Value = '2.55'
$ws.Range("E25").Value = '  -1.94%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").Value = '25.98'
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("E29").Value = '  -3.82%  '
$ws.Range("D30").Value = '35.29'
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("D32").Value = '49.76'
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").Value = '19.92'
$ws.Range("E33").Value = '  +1.17%  '
$ws.Range("D34").Value = '5.38'
$ws.Range("E34").Value = '  -1.95%  '
$ws.Range("D35").Value = '0.0793'
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").Value = '1.97'
$ws.Range("E37").Value = '  -0.96%  '
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").Value = '22.31'
$ws.Range("E41").Value = '  +5.79%  '
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = '119.33'
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '2.002.89'
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("D46").Value = '3.07'
$ws.Range("E46").Value = '  +1.86%  '
$ws.Range("E47").Value = '  -3.36%  '
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").Value = '9.05'
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("E50").Value = '  -2.79%  '
$ws.Range("D51").Value = '57.09'
$ws.Range("E51").Value = '  +3.37%  '